$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template gains a "Mã số thuế" (tax code) column between "Số điện thoại"
# and "Địa chỉ" -> insert a fresh column E, which pushes the old E/F/G
# ("Địa chỉ", "Ngày áp dụng", "Ngày hết hạn") one slot to the right (F/G/H).
# Insert carries over styles/merges/dimension automatically.
$ws.Columns("E:E").Insert()

# New header text for the inserted column.
$ws.Range("E2").Value = "Mã số thuế"

# Column-width bookkeeping: the new column gets its own (narrower) width and
# the old "Địa chỉ" column (now F) gets a little wider.
$ws.Columns("E:E").ColumnWidth = 18.83
$ws.Columns("F:F").ColumnWidth = 29

# Move the cursor/selection onto the new column, matching the refreshed template.
[void]$ws.Columns("E:E").Select()
